$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "Total"
$ws.Range("E12").Formula = "=SUM(E7:E11)"
